$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Wnt2 -> Fzd3 -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt2"
$ws.Range("C2").Value = "Fzd3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01070233333333333
$ws.Range("H2").Value = 0.032107
$ws.Range("I2").Value = 0.004227647500550067
$ws.Range("J2").Value = 0.004227647500550067
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.162165
$ws.Range("N2").Value = 0.486495
$ws.Range("O2").Value = 0.1006291402646046
$ws.Range("P2").Value = 0.1006291402646046
$ws.Range("Q2").Value = 0.001735543885
$ws.Range("R2").Value = 0.015619894965
$ws.Range("S2").Value = 0.0004254245333221576
$ws.Range("T2").Value = 0.0004254245333221577

# Row 3: ECs -> Wnt2 -> Fzd3 -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt2"
$ws.Range("C3").Value = "Fzd3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01070233333333333
$ws.Range("H3").Value = 0.032107
$ws.Range("I3").Value = 0.004227647500550067
$ws.Range("J3").Value = 0.004227647500550067
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.307583
$ws.Range("N3").Value = 0.922749
$ws.Range("O3").Value = 0.1908661724170313
$ws.Range("P3").Value = 0.1908661724170313
$ws.Range("Q3").Value = 0.003291855793666666
$ws.Range("R3").Value = 0.029626702143
$ws.Range("S3").Value = 0.0008069148967584204
$ws.Range("T3").Value = 0.0008069148967584205

# Row 4: ECs -> Wnt2 -> Fzd3 -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt2"
$ws.Range("C4").Value = "Fzd3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01070233333333333
$ws.Range("H4").Value = 0.032107
$ws.Range("I4").Value = 0.004227647500550067
$ws.Range("J4").Value = 0.004227647500550067
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.141763333333333
$ws.Range("N4").Value = 3.42529
$ws.Range("O4").Value = 0.7085046873183641
$ws.Range("P4").Value = 0.7085046873183641
$ws.Range("Q4").Value = 0.01221953178111111
$ws.Range("R4").Value = 0.10997578603
$ws.Range("S4").Value = 0.002995308070469489
$ws.Range("T4").Value = 0.002995308070469489

# Row 5: FAPs -> Wnt2 -> Fzd3 -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt2"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.520808
$ws.Range("H5").Value = 7.562424
$ws.Range("I5").Value = 0.99577235249945
$ws.Range("J5").Value = 0.99577235249945
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.162165
$ws.Range("N5").Value = 0.486495
$ws.Range("O5").Value = 0.1006291402646046
$ws.Range("P5").Value = 0.1006291402646046
$ws.Range("Q5").Value = 0.40878682932
$ws.Range("R5").Value = 3.67908146388
$ws.Range("S5").Value = 0.1002037157312824
$ws.Range("T5").Value = 0.1002037157312824

# Row 6: FAPs -> Wnt2 -> Fzd3 -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt2"
$ws.Range("C6").Value = "Fzd3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.520808
$ws.Range("H6").Value = 7.562424
$ws.Range("I6").Value = 0.99577235249945
$ws.Range("J6").Value = 0.99577235249945
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.307583
$ws.Range("N6").Value = 0.922749
$ws.Range("O6").Value = 0.1908661724170313
$ws.Range("P6").Value = 0.1908661724170313
$ws.Range("Q6").Value = 0.775357687064
$ws.Range("R6").Value = 6.978219183576
$ws.Range("S6").Value = 0.1900592575202729
$ws.Range("T6").Value = 0.1900592575202729

# Row 7: FAPs -> Wnt2 -> Fzd3 -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt2"
$ws.Range("C7").Value = "Fzd3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.520808
$ws.Range("H7").Value = 7.562424
$ws.Range("I7").Value = 0.99577235249945
$ws.Range("J7").Value = 0.99577235249945
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.141763333333333
$ws.Range("N7").Value = 3.42529
$ws.Range("O7").Value = 0.7085046873183641
$ws.Range("P7").Value = 0.7085046873183641
$ws.Range("Q7").Value = 2.878166144773334
$ws.Range("R7").Value = 25.90349530296
$ws.Range("S7").Value = 0.7055093792478947
$ws.Range("T7").Value = 0.7055093792478947
